# "added user by file"
# The sheet's header cell (A1), previously "students", is renamed to
# "Audience Username". The other rows (user name 1..5) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Audience Username"

# Move the selection back onto the edited header cell (the source file no
# longer carries the stale A7 selection that was left over from a prior
# editing session).
[void]$ws.Range("A1").Select()
